$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 567 (shifts existing rows 567-608 down to 568-609),
# matching the weekly update that adds a newer price record ahead of the older ones.
$ws.Rows.Item(567).Insert()

# Populate the newly inserted row 567. It duplicates the data that used to be in
# row 567 (same market/product/quality/prices/origin), but with an updated
# observation date and volume.
$ws.Cells.Item(567, 1).Value = 3
$ws.Cells.Item(567, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(567, 3).Value = 'Coquimbo'
$ws.Cells.Item(567, 4).Value = 45013
$ws.Cells.Item(567, 5).Value = 5
$ws.Cells.Item(567, 6).Value = 'Fruta'
$ws.Cells.Item(567, 7).Value = 100108
$ws.Cells.Item(567, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(567, 9).Value = 100108002
$ws.Cells.Item(567, 10).Value = 'Mango'
$ws.Cells.Item(567, 11).Value = 'Sin especificar'
$ws.Cells.Item(567, 12).Value = 'Primera'
$ws.Cells.Item(567, 13).Value = 228
$ws.Cells.Item(567, 14).Value = 7000
$ws.Cells.Item(567, 15).Value = 7000
$ws.Cells.Item(567, 16).Value = 7000
$ws.Cells.Item(567, 17).Value = '$/bandeja 4 kilos'
$ws.Cells.Item(567, 18).Value = 'Perú'
$ws.Cells.Item(567, 19).Value = 1750
$ws.Cells.Item(567, 20).Value = 4
